$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 44 (shifts existing rows 44-78 down to 45-79)
$ws.Rows.Item(44).Insert()

# Expand Table1 to cover the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A2:X79"))

# Populate the new row with the solved LeetCode problem 572 details
# (value-set order matches the order new shared strings were appended)
$ws.Range("A44").Value = "Binary Tree"
$ws.Range("B44").Value = 572
$ws.Range("C44").Value = "572 - SubTree Of Another Tree"
$ws.Range("D44").Value = "Easy"
$ws.Range("E44").Value = "DFS with nested subtree searches"
$ws.Range("F44").Value = "O(N*M) time, O(N) memory"
$ws.Range("G44").Value = "O(N + M)"
$ws.Range("H44").Value = "Tree Hashing"
$ws.Range("I44").Value = "O(N + M)"
$ws.Range("M44").Value = "45 minutes"
$ws.Range("L44").Value = "Did nto fully solve this one or fully understand the optimal tree hashing solution"

# Highlight the new row red like other "didn't fully solve" rows
$ws.Range("A44").Interior.Color = 255
$ws.Range("B44").Interior.Color = 255
$ws.Range("C44").Interior.Color = 255
$ws.Range("D44").Interior.Color = 255
$ws.Range("E44").Interior.Color = 255
$ws.Range("F44").Interior.Color = 255
$ws.Range("G44").Interior.Color = 255
$ws.Range("H44").Interior.Color = 255
$ws.Range("I44").Interior.Color = 255
$ws.Range("K44").Interior.Color = 255
$ws.Range("L44").Interior.Color = 255
$ws.Range("M44").Interior.Color = 255

# Taller row to fit wrapped notes, matching the other detailed rows
$ws.Rows.Item(44).RowHeight = 58

# Update the view/selection to reflect where the user was working
$win = $excel.ActiveWindow
$win.ScrollRow = 27
[void]$ws.Range("E33").Select()
